# INFY Quarterly Financials update
# - Inserts two new columns (D:E) for the two newest quarterly periods
#   (period-ending dates 2019-01-31 and 2018-10-31), shifting the
#   previously-reported quarters two columns to the right.
# - Refreshes every quarter's reported figures across the Income
#   Statement, Balance Sheet and Cash Flow Statement blocks (rows 7-102)
#   with the latest pulled financial data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new columns before column D (old D:K shifts to F:M)
$ws.Columns("D:E").Insert()

# New columns should carry the same number formatting as the columns
# they were inserted in front of (col F, the old col D)
$ws.Columns("F").Copy()
$ws.Columns("D:E").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$rowData = [ordered]@{
    7 = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
    8 = @(3094400, 2980100, 2765900, 2614800, 2573000, 2540200, 2469500, 2359100, 2380200, 2385300)
    9 = @(2033900, 1929100, 1785500, 1677900, 1663500, 1631200, 1576100, 1484100, 1493800, 1510600)
    10 = @(1060500, 1051000, 980400, 936900, 909500, 909000, 893300, 875000, 886500, 874800)
    11 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    12 = @("NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA", "NA")
    13 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    14 = @(77400, 11000, 49300, 17400, 4200, 6100, 10000, 8000, 5200, 3700)
    15 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    16 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    17 = @(2461200, 2272400, 2148900, 1985200, 1948500, 1926200, 1885300, 1778700, 1783000, 1791500)
    18 = @(633200, 707700, 617000, 629600, 624500, 614000, 584200, 580400, 597200, 593800)
    19 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    20 = @(108900, 106900, 105000, 94300, 139100, 127700, 117700, 99400, 113000, 104300)
    21 = @(826000, 881600, 785000, 790200, 835600, 807600, 767000, 741200, 769900, 756500)
    22 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    23 = @(742100, 814500, 722000, 723900, 763600, 741700, 701900, 679800, 710200, 698100)
    24 = @(220100, 220200, 199700, 190300, 22000, 202900, 198200, 183300, 199300, 201200)
    25 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    26 = @(522000, 594300, 522300, 533600, 741700, 538800, 503600, 496500, 511000, 496900)
    27 = @(521900, 594300, 522300, 533600, 741700, 538800, 503600, 496500, 511000, 496900)
    28 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    29 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    30 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    31 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    32 = @(-108900, -106900, -105000, -94300, -139100, -127700, -117700, -99400, -113000, -104300)
    33 = @(521900, 594300, 522300, 533600, 741700, 538800, 503600, 496500, 511000, 496900)
    34 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    35 = @(521900, 594300, 522300, 533600, 741700, 538800, 503600, 496500, 511000, 496900)
    38 = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
    39 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    40 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    41 = @(1530700, 1678800, 1243700, 1827000, 1635700, 2027100, 2153400, 1972900, 2809300, 3546800)
    42 = @(2216700, 1950800, 2074700, 1888000, 1623400, 3020400, 2604500, 2439900, 1319600, 1050700)
    43 = @(3239000, 3224600, 2998000, 2811000, 3019400, 2965900, 2827700, 2577400, 2721800, 2597400)
    44 = @(10100, 11400, 15000, 17200, 10000, 15900, 12600, 18100, 13400, 9400)
    45 = @(488200, 678600, 678800, 689200, 393900, 391400, 394500, 392300, 355900, 320800)
    46 = @(7484800, 7544200, 7010200, 7232500, 6682400, 8420800, 7992600, 7400500, 7220000, 7525100)
    47 = @(1753000, 1785100, 1846800, 1921700, 1686600, 1800000, 1758300, 1680900, 1497500, 1006500)
    48 = @(1833500, 1789600, 1763000, 1755900, 1711300, 1712900, 1713200, 1614500, 1572300, 1542900)
    49 = @(627900, 413600, 399700, 355400, 625500, 649000, 640700, 610200, 636800, 644200)
    50 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    51 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    52 = @(280800, 299500, 274000, 286600, 277200, 211300, 204200, 180200, 185600, 182000)
    53 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    54 = @(11980000, 11831900, 11293700, 11552100, 10983100, 12793900, 12309100, 11486300, 11112200, 10900800)
    55 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    56 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    57 = @(220500, 172500, 115400, 100400, 72600, 77800, 37600, 50600, 46200, 42300)
    58 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    59 = @(2208000, 2111400, 2275400, 1939200, 1960500, 2204900, 2319000, 1880400, 1982800, 1872400)
    60 = @(2428600, 2284000, 2390800, 2039600, 2033100, 2282700, 2356500, 1931000, 2029000, 1914700)
    61 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    62 = @(147100, 127400, 127700, 124500, 133200, 46600, 45400, 49600, 55100, 53200)
    63 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    64 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    65 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    66 = @(2583400, 2411500, 2518600, 2164200, 2166300, 2329200, 2402000, 1980600, 2084100, 1967900)
    67 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    68 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    69 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    70 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    71 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    72 = @(8894200, 8902300, 8465300, 9092400, 8558700, 9862400, 9323500, 8964700, 8468100, 8374500)
    73 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    74 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    75 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    76 = @(9396500, 9420400, 8775100, 9387900, 8816800, 10464700, 9907100, 9505700, 9028100, 8932900)
    77 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    80 = @(43465, 43373, 43281, 43190, 43100, 43008, 42916, 42825, 42735, 42643)
    81 = @(521900, 594300, 522300, 533600, 741700, 538800, 503600, 496500, 511000, 496900)
    82 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    83 = @(83900, 67100, 63000, 66400, 72000, 65900, 65100, 61500, 59700, 58400)
    84 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    85 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    86 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    87 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    88 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    89 = @(628400, 451900, 616900, 480100, 615600, 409400, 600200, 499500, 510000, 509600)
    90 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    91 = @(-78100, -80100, -77700, -90200, -60000, -58700, -80000, -91400, -86500, -84100)
    92 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    93 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    94 = @(-297400, -83600, -128300, -596300, 1375300, -280800, -48300, -975200, -861100, -320800)
    95 = @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
    96 = @(-529100, -190900, -958600, 700, -492800, -100900, -486300, 0, -417400, -90100)
    97 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    98 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    99 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    100 = @(-528200, -190900, -958600, 1400, -2379200, -100900, -486300, 0, -417400, -90100)
    101 = @(-14500, 15200, -5900, 7800, -6100, 4500, 5500, -5000, -5800, -4700)
    102 = @(-211700, 192600, -475900, -107000, -394500, 32100, 71100, -480600, -774300, 94000)
}


foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $cell = $ws.Cells.Item([int]$r, 4 + $i)
        if ($null -eq $vals[$i]) {
            $cell.Value = $null
        } else {
            $cell.Value = $vals[$i]
        }
    }
}

# Keep the used-range dimension in sync (A5:K102 -> A5:M102)
$ws.Range("A5:M102").Select() | Out-Null
